$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new row 13 data (name, phone number, contribution amount).
# D2:D13 and C14 are formulas (C14/2 and SUM(C2:C13)) that will
# recalculate automatically once C13 is populated.
$ws.Range("A13").Value = "Lalisee Magarsaa"
$ws.Range("B13").Value = 954846351
$ws.Range("C13").Value = 5000

# Match the number format used for the other phone-number cells in column B.
$ws.Range("B13").NumberFormat = $ws.Range("B12").NumberFormat

# Update the active selection to match the saved view state.
$ws.Range("C17").Select()
